$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").NumberFormat = "@"

# Row 2: Santa Catarina
$ws.Range("A2").Value = "Santa Catarina"
$ws.Range("C2").Value = "01/07/2025"
$ws.Range("D2").Value = 55.83
$ws.Range("E2").Value = "1º"

# Row 3: Distrito Federal
$ws.Range("A3").Value = "Distrito Federal"
$ws.Range("C3").Value = "01/07/2025"
$ws.Range("D3").Value = 55.6
$ws.Range("E3").Value = "2º"

# Row 4: São Paulo
$ws.Range("A4").Value = "São Paulo"
$ws.Range("C4").Value = "01/07/2025"
$ws.Range("D4").Value = 55.6
$ws.Range("E4").Value = "3º"

# Row 5: Goiás
$ws.Range("C5").Value = "01/07/2025"
$ws.Range("D5").Value = 54.52
$ws.Range("E5").Value = "4º"

# Row 6: Paraná
$ws.Range("A6").Value = "Paraná"
$ws.Range("C6").Value = "01/07/2025"
$ws.Range("D6").Value = 54.5
$ws.Range("E6").Value = "5º"

# Row 7: Rio Grande do Sul
$ws.Range("A7").Value = "Rio Grande do Sul"
$ws.Range("C7").Value = "01/07/2025"
$ws.Range("D7").Value = 54.19
$ws.Range("E7").Value = "6º"

# Row 8: Sergipe
$ws.Range("C8").Value = "01/07/2025"
$ws.Range("D8").Value = 43.22
$ws.Range("E8").Value = "23º"

# Row 9: Brasil
$ws.Range("C9").Value = "01/07/2025"
$ws.Range("D9").Value = 50.95

# Row 10: Nordeste
$ws.Range("C10").Value = "01/07/2025"
$ws.Range("D10").Value = 44.31
